# Update "想去人数" (F column) figures on the 展览 and 全部类型 sheets
# to reflect the latest generated output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (rows 4-18) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 84
$ws1.Range("F5").Value = 11
$ws1.Range("F6").Value = 12
$ws1.Range("F7").Value = 552
$ws1.Range("F8").Value = 7777
$ws1.Range("F9").Value = 515
$ws1.Range("F10").Value = 213
$ws1.Range("F11").Value = 1086
$ws1.Range("F12").Value = 702
$ws1.Range("F13").Value = 26
$ws1.Range("F15").Value = 186
$ws1.Range("F16").Value = 36
$ws1.Range("F18").Value = 784

# --- Sheet: 全部类型 (rows 4-19) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 84
$ws4.Range("F5").Value = 11
$ws4.Range("F6").Value = 12
$ws4.Range("F8").Value = 552
$ws4.Range("F9").Value = 7777
$ws4.Range("F10").Value = 516
$ws4.Range("F11").Value = 213
$ws4.Range("F12").Value = 1086
$ws4.Range("F13").Value = 702
$ws4.Range("F14").Value = 26
$ws4.Range("F16").Value = 186
$ws4.Range("F17").Value = 36
$ws4.Range("F19").Value = 784
